$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "A new Stock Maintenance System of a Electronics shop is to replace which is very efficient.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A new Stock Maintenance System of a Electronics shop is to replace which is very efficient.",
    2)

$d.Content.Find.Execute(
    "2.1.1 Technology Constraints - Specify the technology stack and any limitations associated with it",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2.1.1 Technology Constraints - Specify the technology stack and any limitations associated with it",
    2)

$d.Content.Find.Execute(
    "Define user roles(admin, manager, stock clerk) with specific permissions,- Implement secure login mechanism.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Define user roles(admin, manager, stock clerk) with specific permissions,- Implement secure login mechanism.",
    2)

$d.Content.Find.Execute(
    "3. Non -Functional Requirements:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3. Non -Functional Requirements",
    2)
